# Update the "Horarios Línea 141" workbook with the latest scrape (01:07:17).
# Sheet "LP1912" gains a new arrival row (111 min, 215_ALUAR) and
# sheet "LP1912-215" gains the same new arrival row; sheet "6203-6173"
# only gets its "Última actualización" timestamp refreshed.

$wb = $excel.ActiveWorkbook

$newTime = "01:07:17"

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 3"

$ws1.Cells.Item(8, 1).Value = $newTime
$ws1.Cells.Item(8, 2).Value = "02:58"
$ws1.Cells.Item(8, 3).Value = "215_ALUAR"
$ws1.Cells.Item(8, 4).Value = 111
$ws1.Cells.Item(8, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Cells.Item(7, 1).Value = $newTime
$ws2.Cells.Item(7, 2).Value = "02:58"
$ws2.Cells.Item(7, 3).Value = "215_ALUAR"
$ws2.Cells.Item(7, 4).Value = 111
$ws2.Cells.Item(7, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
